# "Reading Excel to Two-dim"
# Update the invalidLoginTest sheet (the active sheet): change saul's
# password cell to a numeric value, and add a new "bala" test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: password becomes a plain number instead of the text "saul1223"
$ws.Range("B2").Value = 123344

# New row 4: another invalid-credentials test case
$ws.Range("A4").Value = "bala"
$ws.Range("B4").Value = "bala123"
$ws.Range("C4").Value = "Invalid credentials"

# Move / leave the selection on B2, matching the saved view state
$ws.Range("B2").Select()
